$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 375-385: date-serial, B (nuovi pos.), C (somma mobile 7gg.), D (somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44449, 6, 71, 98.10830604264255),
    @(44450, 24, 81, 111.9263773162542),
    @(44451, 11, 90, 124.3626414625047),
    @(44452, 9, 77, 106.3991488068095),
    @(44453, 2, 74, 102.2537274247261),
    @(44454, 1, 71, 98.10830604264255),
    @(44455, 10, 63, 87.05384902375327),
    @(44456, 2, 59, 81.52662051430862),
    @(44457, 14, 49, 67.70854924069698),
    @(44458, 3, 41, 56.65409222180768),
    @(44459, 15, 47, 64.94493498597465)
)

$startRow = 375
$templateRow = 374

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $serial = $row[0]
    $b = $row[1]
    $c = $row[2]
    $d = $row[3]

    # Copy the formatting from the last existing data row (374, columns A:D) so the
    # new rows inherit the exact same style (borders/font/alignment/number format)
    # rather than creating brand-new style entries.
    $ws.Range("A$templateRow`:D$templateRow").Copy()
    $ws.Range("A$r`:D$r").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($r, 1).Value = $serial
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}
